# Fix a typo in the Korean "press S to start" dialogue string.
# Cell D6 on the DIALOGUE sheet held a copy-paste of the Korean
# "press R to restart" text ('R' 키를 눌러 시작하세요) instead of the
# correct "press S to start" text ('S' 키를 눌러 시작하세요).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Leading "'" keeps the cell's existing quote-prefix text style (it is not
# stored in the cell's text, just forces/kept as literal-text entry).
$ws.Range("D6").Formula = "'" + [char]0x2018 + "S" + [char]0x2019 + " 키를 눌러 시작하세요"

# Restore the active selection to match the post-edit state.
$ws.Range("D15").Select()
